# "added remove member modal (frontend only)" -- mark the four tasks that
# implement the remove-member-modal front-end work as Complete, give them a
# "Date Last Updated" value, and clean up the stray blank Notes cell that
# used to live in D23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Use cells that already carry the styles we want ("Complete" status font in
# B2, and the left-aligned short-date format in D28) as format donors so the
# workbook's existing style table (cellXfs) is reused instead of growing.
$completeStyleSource = $ws.Range("B2")
$dateStyleSource = $ws.Range("D28")

$rows = @(22, 23, 32, 33)
$dates = @("3/4/2024", "3/4/2024", "4/4/2024", "4/4/2024")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]

    $statusCell = $ws.Range("B$r")
    $completeStyleSource.Copy()
    $statusCell.PasteSpecial(-4122)
    $statusCell.Value = "Complete"

    $dateCell = $ws.Range("D$r")
    $dateStyleSource.Copy()
    $dateCell.PasteSpecial(-4122)
    $dateCell.Value = $dates[$i]
}

$excel.CutCopyMode = 0

# Reflect the author's final on-screen selection when the file was saved.
$ws.Range("C29").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
